# Rearranged components in BOM
# Rows 11-22 (resistor section) get reshuffled into a new order.
# Mapping below: key = destination row (after edit), value = source row (before edit),
# both within the same A:F block of the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 11
$lastRow  = 22

# --- snapshot current values + row heights for the affected block --------
$orig = $ws.Range("A" + $firstRow + ":F" + $lastRow).Value2

$origHeights = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $origHeights[$r] = $ws.Rows.Item($r).RowHeight
}

# destination row -> source row (content that should end up there)
$rowMap = @{
    11 = 11
    12 = 20
    13 = 15
    14 = 16
    15 = 12
    16 = 17
    17 = 13
    18 = 22
    19 = 19
    20 = 14
    21 = 18
    22 = 21
}

# --- write the new values back out ---------------------------------------
# NOTE: Range.Value2 hands back a 1-indexed (SAFEARRAY-style) array, but a
# freshly allocated 'New-Object object[,]' is 0-indexed - keep that straight.
$nRows = $orig.GetLength(0)
$nCols = $orig.GetLength(1)
$newVals = New-Object 'object[,]' $nRows, $nCols
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $srcIdx  = $srcRow  - $firstRow + 1
    $destIdx0 = $destRow - $firstRow
    for ($c = 1; $c -le $nCols; $c++) {
        $newVals[$destIdx0, $c - 1] = $orig[$srcIdx, $c]
    }
}
$ws.Range("A" + $firstRow + ":F" + $lastRow).Value2 = $newVals

# --- move the row heights along with their content ------------------------
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $ws.Rows.Item($destRow).RowHeight = $origHeights[$srcRow]
}

# --- hyperlinks: every hyperlinked cell in this sheet has TextToDisplay ---
# --- equal to its target URL, so they can be safely rebuilt from the     ---
# --- (now-reordered) cell text. Deleting any Hyperlinks collection wipes  ---
# --- the whole sheet's hyperlinks in this engine, so delete once and      ---
# --- rebuild every one of them.                                          ---
$hyperlinkRows = @(9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,40)

$ws.Range("E9").Hyperlinks.Delete() | Out-Null

foreach ($r in $hyperlinkRows) {
    $cell = $ws.Cells.Item($r, 5)
    $url = $cell.Value2
    if ($url -ne $null -and $url -ne "") {
        $ws.Hyperlinks.Add($cell, $url) | Out-Null
    }
}

# --- selection: user ends up with row 18 fully selected -------------------
$ws.Rows.Item(18).Select() | Out-Null
